# New entry in ComputerFolders
#
# A new computer ("esmith13laptop") is added to the sheet. Excel represents
# this as a brand-new column inserted right before the old column J
# (bothma-desktop), pushing the existing J/K columns to K/L, and then filling
# in the three rows of data that this computer has values for.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column J; this shifts the former J -> K and
# K -> L (formatting/styles travel with the shifted cells automatically).
$ws.Columns("J").Insert()

# Populate the new computer's entries.
$ws.Range("J1").Value = "esmith13laptop"
$ws.Range("J5").Value = "C:\E\Dropbox\Lab\[07] Transcription\LivemRNAData"
$ws.Range("J8").Value = "C:\E\GitHub\Lab\mRNADynamics"

# Match the column width Excel auto-computed for the new column's content.
$ws.Columns("J").ColumnWidth = 47.666666666666664

# Restore the view/selection state roughly to where the author ended up.
$excel.ActiveWindow.Panes.Item(2).ScrollColumn = 9
$ws.Range("J14").Select() | Out-Null
